$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- 1. Company name update ---
$ws.Range("E7").Value = "RENT INDUSTRIAL SAS"

# --- 2. Header summary numbers ---
$ws.Range("E11").Value = 639106      # VALOR MORA
$ws.Range("C13").Value = 5           # Cant. Trabajadores
$ws.Range("F13").Value = 17          # Cant. Periodos

# --- 3. Insert two extra data rows before the closing/footer row (old row 31) ---
# Row 31 (JOSE FRANCISCO SUAREZ CASTRO / special bottom-border styling) needs to
# stay last among the data rows, so insert the two new rows just above it.
$ws.Rows("31:32").Insert()

# Copy the formatting of a normal interior data row (row 16) onto the two newly
# inserted blank rows so borders/fills/number formats match the rest of the table.
$ws.Range("B16:J16").Copy()
$ws.Range("B31:J31").PasteSpecial(-4122)
$ws.Range("B16:J16").Copy()
$ws.Range("B32:J32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4. Rewrite the full worker/period data block (rows 16-33) ---
$data = @(
    @("CC", "73215869",   "MARCOS ANTONIO SALGADO TORRES", "1911", 5521,  1200000),
    @("CC", "1047458473", "DARWIN GARCIA GUERRERO",         "1911", 5521,  1100000),
    @("CC", "9097751",    "BORIS ESCALANTE THORRENS",       "2308", 52000, 1456000),
    @("CC", "9097751",    "BORIS ESCALANTE THORRENS",       "2307", 52000, 1456000),
    @("CC", "9097751",    "BORIS ESCALANTE THORRENS",       "2306", 31200, 1456000),
    @("CC", "1143344676", "DAYRO MIGUEL CARDALES JULIO",    "2207", 40000, 1000000),
    @("CC", "1143344676", "DAYRO MIGUEL CARDALES JULIO",    "2206", 40000, 1000000),
    @("CC", "1143344676", "DAYRO MIGUEL CARDALES JULIO",    "2205", 40000, 1000000),
    @("CC", "1143344676", "DAYRO MIGUEL CARDALES JULIO",    "2204", 40000, 1000000),
    @("CC", "1143344676", "DAYRO MIGUEL CARDALES JULIO",    "2203", 40000, 1000000),
    @("CC", "1143344676", "DAYRO MIGUEL CARDALES JULIO",    "2202", 40000, 1000000),
    @("CC", "1143344676", "DAYRO MIGUEL CARDALES JULIO",    "2201", 40000, 1000000),
    @("CC", "1143344676", "DAYRO MIGUEL CARDALES JULIO",    "2112", 40000, 1000000),
    @("CC", "1143344676", "DAYRO MIGUEL CARDALES JULIO",    "2111", 40000, 1000000),
    @("CC", "1143344676", "DAYRO MIGUEL CARDALES JULIO",    "2110", 40000, 1000000),
    @("CC", "1143344676", "DAYRO MIGUEL CARDALES JULIO",    "2109", 40000, 1000000),
    @("CC", "1143344676", "DAYRO MIGUEL CARDALES JULIO",    "2108", 37333, 1000000),
    @("CC", "15668088",   "JOSE FRANCISCO SUAREZ CASTRO",   "2401", 15531, 1456000)
)

$r = 16
foreach ($row in $data) {
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $r = $r + 1
}

# --- 5. Column D width ---
$ws.Columns("D:D").ColumnWidth = 34.2
